$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = '28.161.45'
$d.Style = "Normal"
$ws.Range("E2").Value = '  +5.28%  '

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = '1.781.63'
$d.Style = "Normal"
$ws.Range("E3").Value = '  +2.86%  '

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = '1.000'
$d.Style = "Normal"
$ws.Range("E4").Value = '  +0.32%  '

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = '244.45'
$d.Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = '1.000'
$d.Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = '0.4910'
$d.Style = "Normal"
$ws.Range("E7").Value = '  -0.29%  '

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = '0.2678'
$d.Style = "Normal"
$ws.Range("E8").Value = '  +2.07%  '

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = '0.06276'
$d.Style = "Normal"
$ws.Range("E9").Value = '  +0.84%  '

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = '1.782.46'
$d.Style = "Normal"
$ws.Range("E10").Value = '  +3.03%  '

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = '16.46'
$d.Style = "Normal"
$ws.Range("E11").Value = '  +3.69%  '

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = '0.07026'
$d.Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = '0.6274'
$d.Style = "Normal"
$ws.Range("E13").Value = '  +2.46%  '

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = '4.658'
$d.Style = "Normal"
$ws.Range("E14").Value = '  +3.51%  '

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = '79.96'
$d.Style = "Normal"
$ws.Range("E15").Value = '  +3.40%  '

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = '28.135.44'
$d.Style = "Normal"
$ws.Range("E16").Value = '  +6.08%  '

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = '1.001'
$d.Style = "Normal"
$ws.Range("E17").Value = '  +0.25%  '

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = '1.000'
$d.Style = "Normal"
$ws.Range("E18").Value = '  +0.33%  '

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = '0.000007251'
$d.Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = '12.05'
$d.Style = "Normal"
$ws.Range("E20").Value = '  +5.54%  '

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = '2.006.94'
$d.Style = "Normal"
$ws.Range("E21").Value = '  +2.94%  '

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = '4.563'
$d.Style = "Normal"
$ws.Range("E22").Value = '  +1.85%  '

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = '8.737'
$d.Style = "Normal"
$ws.Range("E23").Value = '  +1.99%  '

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = '5.255'
$d.Style = "Normal"
$ws.Range("E24").Value = '  +2.95%  '

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = '141.10'
$d.Style = "Normal"
$ws.Range("E25").Value = '  +2.13%  '

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = '15.76'
$d.Style = "Normal"
$ws.Range("E26").Value = '  +2.69%  '

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = '1.858'
$d.Style = "Normal"
$ws.Range("E27").Value = '  +4.79%  '

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = '109.39'
$d.Style = "Normal"
$ws.Range("E28").Value = '  +2.63%  '

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = '1.386'
$d.Style = "Normal"
$ws.Range("E29").Value = '  -0.22%  '

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = '4.191'
$d.Style = "Normal"
$ws.Range("E30").Value = '  +6.63%  '

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = '0.08261'
$d.Style = "Normal"
$ws.Range("E31").Value = '  +3.36%  '

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = '3.762'
$d.Style = "Normal"
$ws.Range("E32").Value = '  +2.26%  '

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = '0.04890'
$d.Style = "Normal"
$ws.Range("E33").Value = '  +9.13%  '

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = '1.074'
$d.Style = "Normal"
$ws.Range("E34").Value = '  +7.19%  '

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = '2.622'
$d.Style = "Normal"
$ws.Range("E35").Value = '  +0.12%  '

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = '0.6509'
$d.Style = "Normal"
$ws.Range("E36").Value = '  +4.33%  '

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = '0.9510'
$d.Style = "Normal"
$ws.Range("E37").Value = '  +1.79%  '

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = '2.595'
$d.Style = "Normal"
$ws.Range("E38").Value = '  +7.27%  '

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = '2.041'
$d.Style = "Normal"
$ws.Range("E39").Value = '  -0.29%  '

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = '5.891'
$d.Style = "Normal"
$ws.Range("E40").Value = '  +4.91%  '

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = '0.01550'
$d.Style = "Normal"
$ws.Range("E41").Value = '  +2.30%  '

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = '0.9998'
$d.Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = '99.80'
$d.Style = "Normal"
$ws.Range("E43").Value = '  +0.43%  '

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = '0.3981'
$d.Style = "Normal"
$ws.Range("E44").Value = '  +3.13%  '

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = '7.174'
$d.Style = "Normal"
$ws.Range("E45").Value = '  +3.79%  '

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = '0.1213'
$d.Style = "Normal"
$ws.Range("E46").Value = '  +4.47%  '

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = '0.05435'
$d.Style = "Normal"
$ws.Range("E47").Value = '  +0.99%  '

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = '8.018'
$d.Style = "Normal"
$ws.Range("E48").Value = '  +2.22%  '

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = '1.297'
$d.Style = "Normal"
$ws.Range("E49").Value = '  +5.01%  '

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = '30.74'
$d.Style = "Normal"
$ws.Range("E50").Value = '  +1.48%  '

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = '52.91'
$d.Style = "Normal"
$ws.Range("E51").Value = '  +2.23%  '
